$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 415; this shifts the existing rows 415-477
# down to 416-478, preserving all of their values (including the date-format
# style on column D).
$ws.Rows(415).Insert()

# Populate the newly inserted row 415 with the new weekly record.
$ws.Cells.Item(415, 1).Value = 7
$ws.Cells.Item(415, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(415, 3).Value = "Ñuble"
$ws.Cells.Item(415, 4).Value = 45127
$ws.Cells.Item(415, 5).Value = 16
$ws.Cells.Item(415, 6).Value = 100112003
$ws.Cells.Item(415, 7).Value = "Ajo"
$ws.Cells.Item(415, 8).Value = "Chino"
$ws.Cells.Item(415, 9).Value = "Primera"
$ws.Cells.Item(415, 10).Value = 50
$ws.Cells.Item(415, 11).Value = 20000
$ws.Cells.Item(415, 12).Value = 20000
$ws.Cells.Item(415, 13).Value = 20000
$ws.Cells.Item(415, 14).Value = "$/malla 10 kilos"
$ws.Cells.Item(415, 15).Value = "China"
$ws.Cells.Item(415, 16).Value = 2000
$ws.Cells.Item(415, 17).Value = 10
$ws.Cells.Item(415, 18).Value = "Hortaliza"
